# UCD - Proccess RE
# Append a new "Week" block (rows 36-43) to the "Như Phương" timelog sheet,
# by duplicating the previous block (rows 28-35) and adjusting a handful of
# cells, then mark H33 as "Inprogress " instead of "Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Như Phương")
$ws.Activate()

# Copy the previous week's block (rows 28:35) down into rows 36:43 so that
# styles, merged cells, row heights, etc. all come along for the ride.
$src = $ws.Range("A28:I35")
$dst = $ws.Range("A36:I43")
$src.Copy($dst)

# --- Row 36 (Monday) -------------------------------------------------
$ws.Range("D36").Value = "Complete the meeting minutes ; `nReview document "
$ws.Range("G36").Value = 5
$ws.Range("I36").Value = "Meeting_customer 8/10/2019`nMeeting_mentor 6/10/2019"

# --- Row 37 (Tuesday) --------------------------------------------------
$ws.Range("D37").ClearContents()
$ws.Range("I37").ClearContents()

# --- Row 39 (Thursday) -------------------------------------------------
$ws.Range("D39").ClearContents()

# --- Row 41 (Saturday) --------------------------------------------------
$ws.Range("D41").ClearContents()
$ws.Range("H41").Value = "Inprogress "
$ws.Range("I41").ClearContents()

# --- Totals row (row 43) ------------------------------------------------
$ws.Range("C43").Formula = "=SUM(F36:F42)"
$ws.Range("H43").Formula = "=SUM(G36:G42)"

# Existing row 33 status flips from Done to Inprogress.
$ws.Range("H33").Value = "Inprogress "

# Extend the data validations that applied to the previous block so they
# also cover the newly added rows.
$ws.Range("H4:H10,H12:H18,H20:H26,H28:H34,H36:H42").Validation.Delete()
$ws.Range("H4:H10,H12:H18,H20:H26,H28:H34,H36:H42").Validation.Add(3, 1, 1, "Done,Inprogress ")

$ws.Range("C4:C9,C12:C17,C20:C25,C28:C33,C36:C41").Validation.Delete()
$ws.Range("C4:C9,C12:C17,C20:C25,C28:C33,C36:C41").Validation.Add(3, 1, 1, "Project Management, Requirement, Architecture and Desgin, Implementation, Testing, Training, Meetting Customer, Meeting Mentor")

$ws.Range("C10,C18,C26,C34,C42").Validation.Delete()
$ws.Range("C10,C18,C26,C34,C42").Validation.Add(3, 1, 1, "Project Management, Requirement, Architecture and Desgin, Implementation, Testing, Training, Meetting Customer, Meetting Mentor")

# Update the view: scroll down to the new block and select I37 like the
# author left it.
$ws.Range("I37").Select()
$excel.ActiveWindow.ScrollRow = 28

$wb.Save()
